# EDD_calculations_material.xlsx — "includes gauge length calculations"
#
# The underlying change is a single input-cell edit on the "fcc" sheet:
# the 2-theta angle in C14 goes from 30 deg to 13 deg, which cascades
# through the dependent formulas (C15, C16, C18:C24, O5:O33, P5:P33, ...)
# so the live "E" column matches the already-present "E when 2q = 13 deg"
# reference column (T). Excel recalculates all of that automatically.
#
# Alongside that, the extra scatter chart ("Chart 5") that lived on the
# bm_flux sheet (the standalone lineMarker chart) was removed, and the
# selection on the fcc sheet moved from P14 to C18 (the newly-relevant
# gauge-length cell).

$wb = $excel.ActiveWorkbook

# --- fcc sheet: core data edit -------------------------------------------
$fcc = $wb.Worksheets.Item("fcc")
$fcc.Activate()
$fcc.Range("C14").Value = 13

# Move the selection to C18, matching the post-edit cursor position.
$fcc.Range("C18").Select()

# --- bm_flux sheet: drop the extra "Chart 5" scatter chart ---------------
$bmflux = $wb.Worksheets.Item("bm_flux")
$bmflux.Activate()
$charts = $bmflux.ChartObjects()
for ($i = 1; $i -le $charts.Count; $i++) {
    $co = $charts.Item($i)
    if ($co.Name -eq "Chart 5") {
        $co.Delete()
    }
}

# Restore "fcc" as the active sheet/selection (matches the saved file's
# activeTab, which stays on fcc both before and after the edit).
$fcc.Activate()
